# Auto-generated Excel COM-interop script to apply scheduled market-data update
# to the Behemoth_Profits workbook (FFXIV leve-crafting profit tracker).
# For each affected (sheet, row), update currentAveragePrice / NQ / HQ price
# columns (H-N) to the freshly scraped market values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 579.3913  # H17: 804.3077 -> 579.3913
$ws.Cells.Item(17, 10).Value = 579.3913  # J17: 804.3077 -> 579.3913
$ws.Cells.Item(17, 12).Value = 1738.1739  # L17: 2412.9231 -> 1738.1739
$ws.Cells.Item(17, 14).Value = -2074.1739  # N17: -2748.9231 -> -2074.1739
$ws.Cells.Item(51, 8).Value = 131428.58  # H51: 154285.72 -> 131428.58
$ws.Cells.Item(51, 9).Value = 20000  # I51: 0 -> 20000
$ws.Cells.Item(51, 10).Value = 150000  # J51: 154285.72 -> 150000
$ws.Cells.Item(51, 11).Value = 20000  # K51: 0 -> 20000
$ws.Cells.Item(51, 12).Value = 150000  # L51: 154285.72 -> 150000
$ws.Cells.Item(51, 13).Value = -19516  # M51: None -> -19516
$ws.Cells.Item(51, 14).Value = -150968  # N51: -155253.72 -> -150968
$ws.Cells.Item(62, 8).Value = 1595  # H62: 1197.5 -> 1595
$ws.Cells.Item(62, 9).Value = 1595  # I62: 1197.5 -> 1595
$ws.Cells.Item(62, 11).Value = 1595  # K62: 1197.5 -> 1595
$ws.Cells.Item(62, 13).Value = -971  # M62: -573.5 -> -971
$ws.Cells.Item(65, 8).Value = 1595  # H65: 1197.5 -> 1595
$ws.Cells.Item(65, 9).Value = 1595  # I65: 1197.5 -> 1595
$ws.Cells.Item(65, 11).Value = 7975  # K65: 5987.5 -> 7975
$ws.Cells.Item(65, 13).Value = -4855  # M65: -2867.5 -> -4855
$ws.Cells.Item(76, 8).Value = 3450  # H76: 3804.889 -> 3450
$ws.Cells.Item(76, 9).Value = 0  # I76: 3655.5 -> 0
$ws.Cells.Item(76, 10).Value = 3450  # J76: 5000 -> 3450
$ws.Cells.Item(76, 11).Value = 0  # K76: 3655.5 -> 0
$ws.Cells.Item(76, 12).Value = 3450  # L76: 5000 -> 3450
$ws.Cells.Item(76, 13).ClearContents()  # M76: -3340.5 -> (removed)
$ws.Cells.Item(76, 14).Value = -4080  # N76: -5630 -> -4080
$ws.Cells.Item(79, 8).Value = 3450  # H79: 3804.889 -> 3450
$ws.Cells.Item(79, 9).Value = 0  # I79: 3655.5 -> 0
$ws.Cells.Item(79, 10).Value = 3450  # J79: 5000 -> 3450
$ws.Cells.Item(79, 11).Value = 0  # K79: 3655.5 -> 0
$ws.Cells.Item(79, 12).Value = 3450  # L79: 5000 -> 3450
$ws.Cells.Item(79, 13).ClearContents()  # M79: -2563.5 -> (removed)
$ws.Cells.Item(79, 14).Value = -5634  # N79: -7184 -> -5634
$ws.Cells.Item(92, 8).Value = 619.4091  # H92: 645.0952 -> 619.4091
$ws.Cells.Item(92, 9).Value = 604.55554  # I92: 635.41174 -> 604.55554
$ws.Cells.Item(92, 11).Value = 604.55554  # K92: 635.41174 -> 604.55554
$ws.Cells.Item(92, 13).Value = 643.44446  # M92: 612.58826 -> 643.44446
$ws.Cells.Item(132, 8).Value = 2129.8696  # H132: 2179.7556 -> 2129.8696
$ws.Cells.Item(132, 9).Value = 1689.4634  # I132: 1734.575 -> 1689.4634
$ws.Cells.Item(132, 11).Value = 5068.3902  # K132: 5203.725 -> 5068.3902
$ws.Cells.Item(132, 13).Value = -2538.3902  # M132: -2673.725 -> -2538.3902
$ws.Cells.Item(135, 8).Value = 5318.727  # H135: 8020.9287 -> 5318.727
$ws.Cells.Item(135, 9).Value = 765.41174  # I135: 929.3 -> 765.41174
$ws.Cells.Item(135, 10).Value = 20800  # J135: 25750 -> 20800
$ws.Cells.Item(135, 11).Value = 6888.70566  # K135: 8363.699999999999 -> 6888.70566
$ws.Cells.Item(135, 12).Value = 187200  # L135: 231750 -> 187200
$ws.Cells.Item(135, 13).Value = -4353.70566  # M135: -5828.699999999999 -> -4353.70566
$ws.Cells.Item(135, 14).Value = -192270  # N135: -236820 -> -192270
$ws.Cells.Item(138, 8).Value = 2284.1353  # H138: 2286.7837 -> 2284.1353
$ws.Cells.Item(138, 9).Value = 1434.2  # I138: 1444 -> 1434.2
$ws.Cells.Item(138, 11).Value = 4302.6  # K138: 4332 -> 4302.6
$ws.Cells.Item(138, 13).Value = 837.3999999999996  # M138: 808 -> 837.3999999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 2406.6  # H45: 2436 -> 2406.6
$ws.Cells.Item(45, 9).Value = 2039.5625  # I45: 2123.2666 -> 2039.5625
$ws.Cells.Item(45, 10).Value = 3874.75  # J45: 3999.6667 -> 3874.75
$ws.Cells.Item(45, 11).Value = 2039.5625  # K45: 2123.2666 -> 2039.5625
$ws.Cells.Item(45, 12).Value = 3874.75  # L45: 3999.6667 -> 3874.75
$ws.Cells.Item(45, 13).Value = -1662.5625  # M45: -1746.2666 -> -1662.5625
$ws.Cells.Item(45, 14).Value = -4628.75  # N45: -4753.6667 -> -4628.75
$ws.Cells.Item(61, 8).Value = 26375882  # H61: 31320648 -> 26375882
$ws.Cells.Item(61, 9).Value = 38466908  # I61: 50005840 -> 38466908
$ws.Cells.Item(61, 11).Value = 38466908  # K61: 50005840 -> 38466908
$ws.Cells.Item(61, 13).Value = -38466696  # M61: -50005628 -> -38466696
$ws.Cells.Item(122, 8).Value = 1564.8  # H122: 1583.1111 -> 1564.8
$ws.Cells.Item(122, 9).Value = 1457  # I122: 1466.5 -> 1457
$ws.Cells.Item(122, 11).Value = 4371  # K122: 4399.5 -> 4371
$ws.Cells.Item(122, 13).Value = -1921  # M122: -1949.5 -> -1921
$ws.Cells.Item(123, 8).Value = 0  # H123: 88000 -> 0
$ws.Cells.Item(123, 10).Value = 0  # J123: 88000 -> 0
$ws.Cells.Item(123, 12).Value = 0  # L123: 88000 -> 0
$ws.Cells.Item(123, 14).ClearContents()  # N123: -97800 -> (removed)
$ws.Cells.Item(132, 8).Value = 6603.037  # H132: 7027.84 -> 6603.037
$ws.Cells.Item(132, 9).Value = 4224.8335  # I132: 4397.294 -> 4224.8335
$ws.Cells.Item(132, 10).Value = 11359.444  # J132: 12617.75 -> 11359.444
$ws.Cells.Item(132, 11).Value = 12674.5005  # K132: 13191.882 -> 12674.5005
$ws.Cells.Item(132, 12).Value = 34078.33199999999  # L132: 37853.25 -> 34078.33199999999
$ws.Cells.Item(132, 13).Value = -10144.5005  # M132: -10661.882 -> -10144.5005
$ws.Cells.Item(132, 14).Value = -39138.33199999999  # N132: -42913.25 -> -39138.33199999999
$ws.Cells.Item(136, 8).Value = 26375882  # H136: 31320648 -> 26375882
$ws.Cells.Item(136, 9).Value = 38466908  # I136: 50005840 -> 38466908
$ws.Cells.Item(136, 11).Value = 115400724  # K136: 150017520 -> 115400724
$ws.Cells.Item(136, 13).Value = -115398174  # M136: -150014970 -> -115398174

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 12500  # H22: 6322 -> 12500
$ws.Cells.Item(22, 9).Value = 12500  # I22: 6322 -> 12500
$ws.Cells.Item(22, 11).Value = 12500  # K22: 6322 -> 12500
$ws.Cells.Item(22, 13).Value = -12327  # M22: -6149 -> -12327
$ws.Cells.Item(94, 8).Value = 1827.3478  # H94: 1863.2858 -> 1827.3478
$ws.Cells.Item(94, 9).Value = 1832.6818  # I94: 1870.95 -> 1832.6818
$ws.Cells.Item(94, 11).Value = 1832.6818  # K94: 1870.95 -> 1832.6818
$ws.Cells.Item(94, 13).Value = -1381.6818  # M94: -1419.95 -> -1381.6818
$ws.Cells.Item(134, 8).Value = 45969.957  # H134: 47936.086 -> 45969.957
$ws.Cells.Item(134, 9).Value = 1537.2  # I134: 1593.5 -> 1537.2
$ws.Cells.Item(134, 11).Value = 4611.6  # K134: 4780.5 -> 4611.6
$ws.Cells.Item(134, 13).Value = -2076.6  # M134: -2245.5 -> -2076.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 349.5  # H22: 327.85715 -> 349.5
$ws.Cells.Item(22, 9).Value = 349.5  # I22: 327.85715 -> 349.5
$ws.Cells.Item(22, 11).Value = 349.5  # K22: 327.85715 -> 349.5
$ws.Cells.Item(22, 13).Value = 0.5  # M22: 22.14285000000001 -> 0.5
$ws.Cells.Item(31, 8).Value = 490311.06  # H31: 500679.38 -> 490311.06
$ws.Cells.Item(31, 10).Value = 869280.4399999999  # J31: 902598.9 -> 869280.4399999999
$ws.Cells.Item(31, 12).Value = 869280.4399999999  # L31: 902598.9 -> 869280.4399999999
$ws.Cells.Item(31, 14).Value = -869870.4399999999  # N31: -903188.9 -> -869870.4399999999
$ws.Cells.Item(34, 8).Value = 490311.06  # H34: 500679.38 -> 490311.06
$ws.Cells.Item(34, 10).Value = 869280.4399999999  # J34: 902598.9 -> 869280.4399999999
$ws.Cells.Item(34, 12).Value = 869280.4399999999  # L34: 902598.9 -> 869280.4399999999
$ws.Cells.Item(34, 14).Value = -869684.4399999999  # N34: -903002.9 -> -869684.4399999999
$ws.Cells.Item(58, 8).Value = 1228.5  # H58: 1228.25 -> 1228.5
$ws.Cells.Item(58, 9).Value = 966.6667  # I58: 966.3333 -> 966.6667
$ws.Cells.Item(58, 11).Value = 966.6667  # K58: 966.3333 -> 966.6667
$ws.Cells.Item(58, 13).Value = -763.6667  # M58: -763.3333 -> -763.6667
$ws.Cells.Item(62, 8).Value = 3750.7778  # H62: 3965.2856 -> 3750.7778
$ws.Cells.Item(62, 9).Value = 3662.8333  # I62: 3994.25 -> 3662.8333
$ws.Cells.Item(62, 11).Value = 3662.8333  # K62: 3994.25 -> 3662.8333
$ws.Cells.Item(62, 13).Value = -3038.8333  # M62: -3370.25 -> -3038.8333
$ws.Cells.Item(65, 8).Value = 3750.7778  # H65: 3965.2856 -> 3750.7778
$ws.Cells.Item(65, 9).Value = 3662.8333  # I65: 3994.25 -> 3662.8333
$ws.Cells.Item(65, 11).Value = 18314.1665  # K65: 19971.25 -> 18314.1665
$ws.Cells.Item(65, 13).Value = -15194.1665  # M65: -16851.25 -> -15194.1665
$ws.Cells.Item(105, 8).Value = 2559.889  # H105: 2439.9092 -> 2559.889
$ws.Cells.Item(105, 9).Value = 2380.6667  # I105: 2485.5 -> 2380.6667
$ws.Cells.Item(105, 10).Value = 2649.5  # J105: 2413.8572 -> 2649.5
$ws.Cells.Item(105, 11).Value = 2380.6667  # K105: 2485.5 -> 2380.6667
$ws.Cells.Item(105, 12).Value = 2649.5  # L105: 2413.8572 -> 2649.5
$ws.Cells.Item(105, 13).Value = -633.6667000000002  # M105: -738.5 -> -633.6667000000002
$ws.Cells.Item(105, 14).Value = -6143.5  # N105: -5907.8572 -> -6143.5
$ws.Cells.Item(132, 8).Value = 3164.25  # H132: 3375.1538 -> 3164.25
$ws.Cells.Item(132, 9).Value = 3087.7  # I132: 3343.9092 -> 3087.7
$ws.Cells.Item(132, 11).Value = 9263.099999999999  # K132: 10031.7276 -> 9263.099999999999
$ws.Cells.Item(132, 13).Value = -6733.099999999999  # M132: -7501.7276 -> -6733.099999999999
$ws.Cells.Item(136, 8).Value = 1228.5  # H136: 1228.25 -> 1228.5
$ws.Cells.Item(136, 9).Value = 966.6667  # I136: 966.3333 -> 966.6667
$ws.Cells.Item(136, 11).Value = 2900.0001  # K136: 2898.9999 -> 2900.0001
$ws.Cells.Item(136, 13).Value = -350.0001000000002  # M136: -348.9998999999998 -> -350.0001000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(6, 8).Value = 336.33334  # H6: 401.8 -> 336.33334
$ws.Cells.Item(6, 10).Value = 504.5  # J6: 1000 -> 504.5
$ws.Cells.Item(6, 12).Value = 1513.5  # L6: 3000 -> 1513.5
$ws.Cells.Item(6, 14).Value = -1739.5  # N6: -3226 -> -1739.5
$ws.Cells.Item(10, 8).Value = 37.285713  # H10: 38.285713 -> 37.285713
$ws.Cells.Item(10, 9).Value = 26.5  # I10: 27.666666 -> 26.5
$ws.Cells.Item(10, 11).Value = 79.5  # K10: 82.99999800000001 -> 79.5
$ws.Cells.Item(10, 13).Value = 59.5  # M10: 56.00000199999999 -> 59.5
$ws.Cells.Item(19, 8).Value = 400  # H19: 0 -> 400
$ws.Cells.Item(19, 10).Value = 400  # J19: 0 -> 400
$ws.Cells.Item(19, 12).Value = 1200  # L19: 0 -> 1200
$ws.Cells.Item(19, 14).Value = -1548  # N19: None -> -1548
$ws.Cells.Item(113, 8).Value = 1911  # H113: 1606.8334 -> 1911
$ws.Cells.Item(113, 10).Value = 2822  # J113: 1910.25 -> 2822
$ws.Cells.Item(113, 12).Value = 8466  # L113: 5730.75 -> 8466
$ws.Cells.Item(113, 14).Value = -12806  # N113: -10070.75 -> -12806
$ws.Cells.Item(114, 8).Value = 1184  # H114: 2305.3333 -> 1184
$ws.Cells.Item(114, 9).Value = 968.4  # I114: 2318.3333 -> 968.4
$ws.Cells.Item(114, 10).Value = 1399.6  # J114: 2266.3333 -> 1399.6
$ws.Cells.Item(114, 11).Value = 2905.2  # K114: 6954.999899999999 -> 2905.2
$ws.Cells.Item(114, 12).Value = 4198.799999999999  # L114: 6798.999899999999 -> 4198.799999999999
$ws.Cells.Item(114, 13).Value = 348.8000000000002  # M114: -3700.999899999999 -> 348.8000000000002
$ws.Cells.Item(114, 14).Value = -10706.8  # N114: -13306.9999 -> -10706.8
$ws.Cells.Item(121, 8).Value = 3541.4285  # H121: 4048.3333 -> 3541.4285
$ws.Cells.Item(121, 10).Value = 5747.5  # J121: 7496.6665 -> 5747.5
$ws.Cells.Item(121, 12).Value = 17242.5  # L121: 22489.9995 -> 17242.5
$ws.Cells.Item(121, 14).Value = -19862.5  # N121: -25109.9995 -> -19862.5
$ws.Cells.Item(131, 8).Value = 5059.65  # H131: 5350.4736 -> 5059.65
$ws.Cells.Item(131, 9).Value = 6514.8  # I131: 7268.6665 -> 6514.8
$ws.Cells.Item(131, 10).Value = 3604.5  # J131: 3624.1 -> 3604.5
$ws.Cells.Item(131, 11).Value = 19544.4  # K131: 21805.9995 -> 19544.4
$ws.Cells.Item(131, 12).Value = 10813.5  # L131: 10872.3 -> 10813.5
$ws.Cells.Item(131, 13).Value = -14504.4  # M131: -16765.9995 -> -14504.4
$ws.Cells.Item(131, 14).Value = -20893.5  # N131: -20952.3 -> -20893.5
$ws.Cells.Item(137, 8).Value = 7671.2856  # H137: 5138.077 -> 7671.2856
$ws.Cells.Item(137, 9).Value = 6559.8  # I137: 5549.6665 -> 6559.8
$ws.Cells.Item(137, 10).Value = 10450  # J137: 4785.2856 -> 10450
$ws.Cells.Item(137, 11).Value = 19679.4  # K137: 16648.9995 -> 19679.4
$ws.Cells.Item(137, 12).Value = 31350  # L137: 14355.8568 -> 31350
$ws.Cells.Item(137, 13).Value = -14579.4  # M137: -11548.9995 -> -14579.4
$ws.Cells.Item(137, 14).Value = -41550  # N137: -24555.8568 -> -41550
$ws.Cells.Item(138, 8).Value = 2599.6  # H138: 2799.6 -> 2599.6
$ws.Cells.Item(138, 9).Value = 1666  # I138: 1999.3334 -> 1666
$ws.Cells.Item(138, 11).Value = 4998  # K138: 5998.0002 -> 4998
$ws.Cells.Item(138, 13).Value = 142  # M138: -858.0002000000004 -> 142

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 217549.86  # H3: 608.8333 -> 217549.86
$ws.Cells.Item(3, 9).Value = 173450  # I3: 251 -> 173450
$ws.Cells.Item(3, 10).Value = 250624.75  # J3: 966.6667 -> 250624.75
$ws.Cells.Item(3, 11).Value = 173450  # K3: 251 -> 173450
$ws.Cells.Item(3, 12).Value = 250624.75  # L3: 966.6667 -> 250624.75
$ws.Cells.Item(3, 13).Value = -173334  # M3: -135 -> -173334
$ws.Cells.Item(3, 14).Value = -250856.75  # N3: -1198.6667 -> -250856.75
$ws.Cells.Item(111, 8).Value = 53719.332  # H111: 54435.668 -> 53719.332
$ws.Cells.Item(111, 10).Value = 53719.332  # J111: 54435.668 -> 53719.332
$ws.Cells.Item(111, 12).Value = 53719.332  # L111: 54435.668 -> 53719.332
$ws.Cells.Item(111, 14).Value = -59853.332  # N111: -60569.668 -> -59853.332
$ws.Cells.Item(126, 8).Value = 5249.5  # H126: 4603.6665 -> 5249.5
$ws.Cells.Item(126, 9).Value = 4500  # I126: 3906 -> 4500
$ws.Cells.Item(126, 11).Value = 13500  # K126: 11718 -> 13500
$ws.Cells.Item(126, 13).Value = -11030  # M126: -9248 -> -11030

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2104.5386  # H16: 2199.6428 -> 2104.5386
$ws.Cells.Item(16, 9).Value = 1787.6364  # I16: 1925 -> 1787.6364
$ws.Cells.Item(16, 11).Value = 1787.6364  # K16: 1925 -> 1787.6364
$ws.Cells.Item(16, 13).Value = -1617.6364  # M16: -1755 -> -1617.6364
$ws.Cells.Item(22, 8).Value = 2918.1365  # H22: 2867.348 -> 2918.1365
$ws.Cells.Item(22, 9).Value = 2960.7856  # I22: 2880.0667 -> 2960.7856
$ws.Cells.Item(22, 11).Value = 2960.7856  # K22: 2880.0667 -> 2960.7856
$ws.Cells.Item(22, 13).Value = -2665.7856  # M22: -2585.0667 -> -2665.7856
$ws.Cells.Item(27, 8).Value = 2918.1365  # H27: 2867.348 -> 2918.1365
$ws.Cells.Item(27, 9).Value = 2960.7856  # I27: 2880.0667 -> 2960.7856
$ws.Cells.Item(27, 11).Value = 2960.7856  # K27: 2880.0667 -> 2960.7856
$ws.Cells.Item(27, 13).Value = -2853.7856  # M27: -2773.0667 -> -2853.7856
$ws.Cells.Item(46, 8).Value = 4470.2856  # H46: 4353.4546 -> 4470.2856
$ws.Cells.Item(46, 10).Value = 4905.5557  # J46: 4605 -> 4905.5557
$ws.Cells.Item(46, 12).Value = 4905.5557  # L46: 4605 -> 4905.5557
$ws.Cells.Item(46, 14).Value = -5281.5557  # N46: -4981 -> -5281.5557
$ws.Cells.Item(116, 8).Value = 129650  # H116: 264666.66 -> 129650
$ws.Cells.Item(116, 10).Value = 129650  # J116: 264666.66 -> 129650
$ws.Cells.Item(116, 12).Value = 129650  # L116: 264666.66 -> 129650
$ws.Cells.Item(116, 14).Value = -138828  # N116: -273844.66 -> -138828
$ws.Cells.Item(136, 8).Value = 91172.16  # H136: 98077.664 -> 91172.16
$ws.Cells.Item(136, 9).Value = 2462  # I136: 2365.875 -> 2462
$ws.Cells.Item(136, 10).Value = 194667.33  # J136: 289501.25 -> 194667.33
$ws.Cells.Item(136, 11).Value = 7386  # K136: 7097.625 -> 7386
$ws.Cells.Item(136, 12).Value = 584001.99  # L136: 868503.75 -> 584001.99
$ws.Cells.Item(136, 13).Value = -4836  # M136: -4547.625 -> -4836
$ws.Cells.Item(136, 14).Value = -589101.99  # N136: -873603.75 -> -589101.99

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4451780.5  # H62: 4658747 -> 4451780.5
$ws.Cells.Item(62, 9).Value = 7460.2163  # I62: 7611.8887 -> 7460.2163
$ws.Cells.Item(62, 10).Value = 25006762  # J62: 28578872 -> 25006762
$ws.Cells.Item(62, 11).Value = 7460.2163  # K62: 7611.8887 -> 7460.2163
$ws.Cells.Item(62, 12).Value = 25006762  # L62: 28578872 -> 25006762
$ws.Cells.Item(62, 13).Value = -6836.2163  # M62: -6987.8887 -> -6836.2163
$ws.Cells.Item(62, 14).Value = -25008010  # N62: -28580120 -> -25008010
$ws.Cells.Item(65, 8).Value = 4451780.5  # H65: 4658747 -> 4451780.5
$ws.Cells.Item(65, 9).Value = 7460.2163  # I65: 7611.8887 -> 7460.2163
$ws.Cells.Item(65, 10).Value = 25006762  # J65: 28578872 -> 25006762
$ws.Cells.Item(65, 11).Value = 37301.0815  # K65: 38059.4435 -> 37301.0815
$ws.Cells.Item(65, 12).Value = 125033810  # L65: 142894360 -> 125033810
$ws.Cells.Item(65, 13).Value = -34181.0815  # M65: -34939.4435 -> -34181.0815
$ws.Cells.Item(65, 14).Value = -125040050  # N65: -142900600 -> -125040050
$ws.Cells.Item(81, 8).Value = 600.2857  # H81: 605.6667 -> 600.2857
$ws.Cells.Item(81, 9).Value = 614.25  # I81: 629.6667 -> 614.25
$ws.Cells.Item(81, 11).Value = 1228.5  # K81: 1259.3334 -> 1228.5
$ws.Cells.Item(81, 13).Value = -167.5  # M81: -198.3334 -> -167.5
$ws.Cells.Item(84, 8).Value = 600.2857  # H84: 605.6667 -> 600.2857
$ws.Cells.Item(84, 9).Value = 614.25  # I84: 629.6667 -> 614.25
$ws.Cells.Item(84, 11).Value = 6142.5  # K84: 6296.666999999999 -> 6142.5
$ws.Cells.Item(84, 13).Value = -838.5  # M84: -992.6669999999995 -> -838.5
$ws.Cells.Item(122, 8).Value = 4136.6343  # H122: 4440.1353 -> 4136.6343
$ws.Cells.Item(122, 9).Value = 1859.04  # I122: 1939.0454 -> 1859.04
$ws.Cells.Item(122, 10).Value = 7695.375  # J122: 8108.4 -> 7695.375
$ws.Cells.Item(122, 11).Value = 5577.12  # K122: 5817.1362 -> 5577.12
$ws.Cells.Item(122, 12).Value = 23086.125  # L122: 24325.2 -> 23086.125
$ws.Cells.Item(122, 13).Value = -3127.12  # M122: -3367.1362 -> -3127.12
$ws.Cells.Item(122, 14).Value = -27986.125  # N122: -29225.2 -> -27986.125
